# Update the "Förändrad" date column (C) for every data row from
# 2023-09-08 (serial 45177) to 2023-09-09 (serial 45178).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C532").Value = 45178
